# SW1116_noCTRL_meas.xlsx bug-fix edit
#
# The source workbook accumulated a tail of stray rows on Sheet1 (rows
# 45-87 only ever held a leftover index number in column A - debris from
# an earlier autofill) that don't belong to the real 44-row dataset the
# other two sheets use. Clean that up, then leave the file the way the
# author left it open: positioned on Sheet1, scrolled down a bit with
# cell E58 selected (rather than back on Sheet3 with the old B2:N44 /
# A2:N44 selections).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet3 = $wb.Worksheets.Item("Sheet3")

# Drop the stray rows 45:87 on Sheet1 - only column A had data there, and
# it was left over past the real 44-row table (dimension collapses back
# to A1:N44 automatically once the rows are gone).
$sheet1.Rows("45:87").Delete()

# Sheet3 was the active/selected tab before; move the active tab back to
# Sheet1 (this also clears Sheet3's tabSelected flag).
$sheet1.Activate()

# Park the selection on E58, as the workbook was left scrolled to that
# area of the (now-shorter) sheet.
$sheet1.Range("E58").Select()
